$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22
$ws.Range("B22").Value = 0.5454545454545454
$ws.Range("D22").Value = 0.5714285714285713

# Row 23
$ws.Range("B23").Value = 0.6923076923076923
$ws.Range("C23").Value = 0.6428571428571429
$ws.Range("D23").Value = 0.6666666666666666

# Row 24
$ws.Range("B24").Value = 0.625
$ws.Range("C24").Value = 0.625
$ws.Range("D24").Value = 0.625
$ws.Range("E24").Value = 0.625

# Row 25
$ws.Range("B25").Value = 0.6188811188811189
$ws.Range("C25").Value = 0.6214285714285714
$ws.Range("D25").Value = 0.619047619047619

# Row 26
$ws.Range("B26").Value = 0.6311188811188811
$ws.Range("C26").Value = 0.625
$ws.Range("D26").Value = 0.6269841269841269
